# The document contains (as plain text, across several differently-
# formatted runs):
#   ...mays plustost il les fault <<corr><del>les</del></corr> bien...
# The first "<" right after "mays plustost il les fault " lives in its
# own run (Courier New, blue 0000ff, 9pt) and must be removed entirely,
# while the text/run that follows ("<corr><del>") must be left untouched
# with its own (different) formatting. A plain text Find/Replace would
# merge runs and lose the surrounding formatting, so instead we locate
# the unique anchor text and delete just that single character's Range,
# which removes the whole (one-character) run cleanly.

$d = $word.ActiveDocument

$anchor = "mays plustost il les fault "
$text = $d.Content.Text
$anchorPos = $text.IndexOf($anchor)

if ($anchorPos -lt 0) {
    throw "Anchor text not found: $anchor"
}

$charStart = $anchorPos + $anchor.Length
$charEnd = $charStart + 1

$target = $d.Range($charStart, $charEnd)

if ($target.Text -ne "<") {
    throw "Unexpected character at target position: [$($target.Text)]"
}

$target.Delete()
